# Informe de Avance.docx -- apply commit:
#   "Incluye agregarCarrito, eliminarCarrito, confirmarCarritoEnvio,
#    confirmarCarritoRetiroLocal, etc."
#
# The headline change is a new bullet under Julian Castro's "Detalle
# participacion" section listing the Comercio cart-related methods. The
# surrounding diff also shows a handful of proof-reading/run-shape
# touch-ups (re-run spell check splits runs differently, a couple of
# redundant grammar-check wrapper runs collapse back into their
# neighbours, and a stray lastRenderedPageBreak hint disappears). We
# replicate all of the textual/structural edits; none of them change
# what the document displays, only how the runs are segmented.

function Get-ParagraphByText($doc, $needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    throw "No paragraph found containing '$needle'"
}

# Swaps the paragraph that contains $needle for the literal OOXML in
# $xmlFragment (one or more sibling <w:p> elements). Because the found
# paragraph's Range is collapsed to its own start before InsertXML is
# called, the paragraph mark that used to end the old paragraph becomes
# the end of the *last* inserted <w:p>, so the old paragraph is fully
# replaced (not just prefixed).
function Replace-ParagraphXml($doc, $needle, $xmlFragment) {
    $p = Get-ParagraphByText $doc $needle
    $r = $p.Range
    $r.Collapse(1)
    $r.InsertXML($xmlFragment)
}

$d = $word.ActiveDocument
$xmlns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "Gianluca " run re-split by the spell checker ----------------
$gianlucaXml = @"
<w:p $xmlns w:rsidR="00980597" w:rsidRDefault="000517B8"><w:pPr><w:pStyle w:val="Ttulo1"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:bookmarkStart w:id="8" w:name="_eej9qx7ryk3t" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="8"/><w:r><w:rPr><w:b w:val="0"/><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">&lt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="003F75B1"><w:rPr><w:b w:val="0"/><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Gianluca</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="003F75B1"><w:rPr><w:b w:val="0"/><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="003F75B1"><w:rPr><w:b w:val="0"/><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Cambareri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b w:val="0"/><w:color w:val="434343"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> &gt;</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "Gianluca" $gianlucaXml

# --- 2. "Detalle participacion (...)" gramStart/End wrapper merges away
$detalleXml = @"
<w:p $xmlns w:rsidR="00980597" w:rsidRPr="00AB4755" w:rsidRDefault="000517B8"><w:pPr><w:pStyle w:val="Ttulo1"/><w:widowControl w:val="0"/><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:bookmarkStart w:id="12" w:name="_6l209sievu4g" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="12"/><w:r w:rsidRPr="00AB4755"><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Detalle participaci&#243;n (por cada integrante) :</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "Detalle participaci" $detalleXml

# --- 3. New bullet "Comercio: agregarCarrito(), ... etc." inserted ---
#        right before the existing "Carrito: agregarItem() ..." bullet.
$carritoNewXml = @"
<w:p $xmlns w:rsidR="00AB4755" w:rsidRDefault="00AB4755" w:rsidP="00AB4755"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Comercio: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>agregarCarrito</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>eliminarCarrito</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>confirmarCarritoEnvio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>confirmarCarritoRetiroLocal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>mostrarCarrito</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>(), etc.</w:t></w:r></w:p>
"@
$carritoOldXml = @"
<w:p $xmlns w:rsidR="00AB4755" w:rsidRDefault="00AB4755" w:rsidP="00AB4755"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Carrito: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>agregarItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>eliminarItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>toString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "agregarItem" ($carritoNewXml + $carritoOldXml)

# --- 4. Kevin Canepa's three "Comercio: traer...()" bullets: the
#        "(" / ")" pair (straddling a gramStart/gramEnd wrapper) merges
#        back into a single "()" run.
$traerArticuloXml = @"
<w:p $xmlns w:rsidR="00E61CCB" w:rsidRDefault="00E61CCB" w:rsidP="00E61CCB"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Comercio: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>traerArticulo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "traerArticulo" $traerArticuloXml

$traerAgendaXml = @"
<w:p $xmlns w:rsidR="00E61CCB" w:rsidRDefault="00E61CCB" w:rsidP="00E61CCB"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Comercio: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>traerAgenda</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "traerAgenda" $traerAgendaXml

$traerTurnosLibresXml = @"
<w:p $xmlns w:rsidR="00E61CCB" w:rsidRDefault="00E61CCB" w:rsidP="00E61CCB"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Comercio: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>traerTurnosLibres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "traerTurnosLibres" $traerTurnosLibresXml

# --- 5. "Articulo: " also gets spell-checked/re-split, and its "()"
#        pair collapses the same way as the bullets above.
$validarCodBarrasXml = @"
<w:p $xmlns w:rsidR="00E61CCB" w:rsidRDefault="00E61CCB" w:rsidP="00E61CCB"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Articulo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>validarCodBarras</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>()</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "validarCodBarras" $validarCodBarrasXml

# --- 6. "Clases: Comercio, ..." -- "Comercio" gets its own spell-check
#        wrapped run instead of being glued to the leading ": ".
$clasesXml = @"
<w:p $xmlns w:rsidR="00980597" w:rsidRPr="000517B8" w:rsidRDefault="000517B8" w:rsidP="000517B8"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Clases</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Comercio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DiaRetiro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Turno</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@
Replace-ParagraphXml $d ": Comercio, " $clasesXml

# --- 7. Drop the stale lastRenderedPageBreak hint before "Inconvenientes"
$inconvenientesXml = @"
<w:p $xmlns w:rsidR="00980597" w:rsidRPr="003F75B1" w:rsidRDefault="000517B8"><w:pPr><w:pStyle w:val="Ttulo1"/><w:widowControl w:val="0"/></w:pPr><w:bookmarkStart w:id="18" w:name="_t04gojr7v82u" w:colFirst="0" w:colLast="0"/><w:bookmarkEnd w:id="18"/><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003F75B1"><w:t>Inconvenientes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="003F75B1"><w:t>:</w:t></w:r></w:p>
"@
Replace-ParagraphXml $d "Inconvenientes" $inconvenientesXml

Write-Output "done"
